# Add a new worksheet "ODI Batting Extra" after the existing "ODI Batting" sheet
# and populate it with per-match batting extras data.

$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Insert the new sheet right after "ODI Batting" so it becomes the 3rd / last tab.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBatting)
$ws.Name = "ODI Batting Extra"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Copy the header formatting (bold, centered, bordered) used by the other sheets.
$odiBatting.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
# Row 2 - match 4231
$ws.Range("A2").Value = "'4231"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'2.96%"
$ws.Range("F2").Value = "NO"

# Row 3 - match 4232
$ws.Range("A3").Value = "'4232"
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'1"
$ws.Range("E3").Value = "'6.90%"
$ws.Range("F3").Value = "NO"

# Row 4 - match 4233 (only match code and man-of-match known)
$ws.Range("A4").Value = "'4233"
$ws.Range("F4").Value = "NO"
